$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 142; A = "Transmission correct"; B = "e0y4cjuz"; C = "Training phase"; D = 1; E = "['Purple', 'Green']"; F = "[['Red', 'Blue'], ['Blue', 'Yellow']]"; G = "[None, None]"; H = "['8', '8']"; I = $null },
    @{ Row = 143; A = "Transmission correct"; B = "e0y4cjuz"; C = "Training phase"; D = 2; E = "['Green', 'Purple']"; F = "[['Red', ''], ['Blue', '']]"; G = "[None, None]"; H = "['2', '5']"; I = "0.23" },
    @{ Row = 144; A = "Transmission M&M"; B = "c2lm5k76"; C = "Training phase"; D = 1; E = "['Purple', 'Green']"; F = "[['Red', ''], ['Blue', '']]"; G = "[None, None]"; H = "['2', '5']"; I = "0.07" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    if ($r.I -ne $null) {
        $ws.Cells.Item($row, 9).NumberFormat = "@"
        $ws.Cells.Item($row, 9).Value = $r.I
        $ws.Cells.Item($row, 9).Style = "Normal"
    }
}
